$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = $true
$ws.Range("B3").Value = "hyper_heuristic"
$ws.Range("C3").Value = 3716.351584315983

$ws.Range("A4").Value = $false
$ws.Range("B4").Value = "hyper_heuristic"
$ws.Range("C4").Value = 3850.610954794818

$ws.Range("A5").Value = $true
$ws.Range("B5").Value = "simple"
$ws.Range("C5").Value = 3064.831742360079

$ws.Range("A6").Value = $false
$ws.Range("B6").Value = "simple"
$ws.Range("C6").Value = 2998.187948059471
